# Project "Sample Project" save: update rule R40's "From" value (cell B11
# on the "Rules" sheet) from the text "R40" to the text "1".
#
# A plain `Range.Value = "1"` assignment would be auto-coerced to the
# *number* 1 by Excel's normal type inference (since "1" parses cleanly as
# a number), which would change the cell's stored type away from a shared
# text string. To keep the cell holding literal text "1" (matching how the
# workbook was actually edited), we round-trip the new text through a
# formula and then paste it back as a value: this keeps the cell's
# underlying type as text instead of letting it fall back to a number,
# while leaving every other attribute of the cell (its style/format) as it
# was.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("B11")

$cell.Formula = '="1"'
$cell.Copy()
$cell.PasteSpecial(-4163)  # xlPasteValues - keep the result as a literal value, not a formula
$excel.CutCopyMode = $false
